$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new column D
$ws.Range("D1").Value = "Correct_answer"

# Rows 2-51 correspond to the "Purple" category -> correct answer "l"
$ws.Range("D2:D51").Value = "l"

# Rows 52-101 correspond to the "Blue" category -> correct answer "s"
$ws.Range("D52:D101").Value = "s"

# Update the view: drop the previous scroll position / selection and
# select cell F5 instead.
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("F5").Select()
